$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 9).Value = 0.3800727954645477
$ws.Cells.Item(2, 10).Value = 0.47906870267432
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.07487166666666667
$ws.Cells.Item(2, 14).Value = 0.224615
$ws.Cells.Item(2, 15).Value = 0.01287435003490057
$ws.Cells.Item(2, 16).Value = 0.01655871537719798
$ws.Cells.Item(2, 17).Value = 0.03467918335277778
$ws.Cells.Item(2, 18).Value = 0.312112650175
$ws.Cells.Item(2, 19).Value = 0.004893190207553757
$ws.Cells.Item(2, 20).Value = 0.007932762293707547
$ws.Cells.Item(3, 9).Value = 0.3800727954645477
$ws.Cells.Item(3, 10).Value = 0.47906870267432
$ws.Cells.Item(3, 15).Value = 0.02600892111095355
$ws.Cells.Item(3, 16).Value = 0.03345212152666174
$ws.Cells.Item(3, 17).Value = 0.07005931496111112
$ws.Cells.Item(3, 18).Value = 0.63053383465
$ws.Cells.Item(3, 19).Value = 0.009885283353657005
$ws.Cells.Item(3, 20).Value = 0.01602586446148153
$ws.Cells.Item(4, 9).Value = 0.3800727954645477
$ws.Cells.Item(4, 10).Value = 0.47906870267432
$ws.Cells.Item(4, 13).Value = 0.8000470000000001
$ws.Cells.Item(4, 14).Value = 2.400141
$ws.Cells.Item(4, 15).Value = 0.1375698656239178
$ws.Cells.Item(4, 16).Value = 0.1769394371887155
$ws.Cells.Item(4, 17).Value = 0.3705671028716667
$ws.Cells.Item(4, 18).Value = 3.335103925845
$ws.Cells.Item(4, 19).Value = 0.05228656339936461
$ws.Cells.Item(4, 20).Value = 0.08476614662592225
$ws.Cells.Item(5, 9).Value = 0.3800727954645477
$ws.Cells.Item(5, 10).Value = 0.47906870267432
$ws.Cells.Item(5, 13).Value = 3.881946
$ws.Cells.Item(5, 14).Value = 7.763892
$ws.Cells.Item(5, 15).Value = 0.6675092708044715
$ws.Cells.Item(5, 16).Value = 0.5723574910282232
$ws.Cells.Item(5, 17).Value = 1.79804621819
$ws.Cells.Item(5, 18).Value = 10.78827730914
$ws.Cells.Item(5, 19).Value = 0.2537021145531573
$ws.Cells.Item(5, 20).Value = 0.2741985606928196
$ws.Cells.Item(6, 9).Value = 0.3800727954645477
$ws.Cells.Item(6, 10).Value = 0.47906870267432
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.9074473333333333
$ws.Cells.Item(6, 14).Value = 2.722342
$ws.Cells.Item(6, 15).Value = 0.1560375924257564
$ws.Cells.Item(6, 16).Value = 0.2006922348792017
$ws.Cells.Item(6, 17).Value = 0.4203129682655555
$ws.Cells.Item(6, 18).Value = 3.78281671439
$ws.Cells.Item(6, 19).Value = 0.05930564395081498
$ws.Cells.Item(6, 20).Value = 0.09614536860038907
$ws.Cells.Item(7, 7).Value = 0.755484
$ws.Cells.Item(7, 8).Value = 1.510968
$ws.Cells.Item(7, 9).Value = 0.6199272045354524
$ws.Cells.Item(7, 10).Value = 0.52093129732568
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.07487166666666667
$ws.Cells.Item(7, 14).Value = 0.224615
$ws.Cells.Item(7, 15).Value = 0.01287435003490057
$ws.Cells.Item(7, 16).Value = 0.01655871537719798
$ws.Cells.Item(7, 17).Value = 0.05656434622000001
$ws.Cells.Item(7, 18).Value = 0.33938607732
$ws.Cells.Item(7, 19).Value = 0.007981159827346814
$ws.Cells.Item(7, 20).Value = 0.008625953083490429
$ws.Cells.Item(8, 7).Value = 0.755484
$ws.Cells.Item(8, 8).Value = 1.510968
$ws.Cells.Item(8, 9).Value = 0.6199272045354524
$ws.Cells.Item(8, 10).Value = 0.52093129732568
$ws.Cells.Item(8, 15).Value = 0.02600892111095355
$ws.Cells.Item(8, 16).Value = 0.03345212152666174
$ws.Cells.Item(8, 17).Value = 0.11427199156
$ws.Cells.Item(8, 18).Value = 0.68563194936
$ws.Cells.Item(8, 19).Value = 0.01612363775729654
$ws.Cells.Item(8, 20).Value = 0.01742625706518021
$ws.Cells.Item(9, 7).Value = 0.755484
$ws.Cells.Item(9, 8).Value = 1.510968
$ws.Cells.Item(9, 9).Value = 0.6199272045354524
$ws.Cells.Item(9, 10).Value = 0.52093129732568
$ws.Cells.Item(9, 13).Value = 0.8000470000000001
$ws.Cells.Item(9, 14).Value = 2.400141
$ws.Cells.Item(9, 15).Value = 0.1375698656239178
$ws.Cells.Item(9, 16).Value = 0.1769394371887155
$ws.Cells.Item(9, 17).Value = 0.6044227077480001
$ws.Cells.Item(9, 18).Value = 3.626536246488
$ws.Cells.Item(9, 19).Value = 0.08528330222455316
$ws.Cells.Item(9, 20).Value = 0.09217329056279322
$ws.Cells.Item(10, 7).Value = 0.755484
$ws.Cells.Item(10, 8).Value = 1.510968
$ws.Cells.Item(10, 9).Value = 0.6199272045354524
$ws.Cells.Item(10, 10).Value = 0.52093129732568
$ws.Cells.Item(10, 13).Value = 3.881946
$ws.Cells.Item(10, 14).Value = 7.763892
$ws.Cells.Item(10, 15).Value = 0.6675092708044715
$ws.Cells.Item(10, 16).Value = 0.5723574910282232
$ws.Cells.Item(10, 17).Value = 2.932748091864
$ws.Cells.Item(10, 18).Value = 11.730992367456
$ws.Cells.Item(10, 19).Value = 0.4138071562513143
$ws.Cells.Item(10, 20).Value = 0.2981589303354036
$ws.Cells.Item(11, 7).Value = 0.755484
$ws.Cells.Item(11, 8).Value = 1.510968
$ws.Cells.Item(11, 9).Value = 0.6199272045354524
$ws.Cells.Item(11, 10).Value = 0.52093129732568
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.9074473333333333
$ws.Cells.Item(11, 14).Value = 2.722342
$ws.Cells.Item(11, 15).Value = 0.1560375924257564
$ws.Cells.Item(11, 16).Value = 0.2006922348792017
$ws.Cells.Item(11, 17).Value = 0.685561941176
$ws.Cells.Item(11, 18).Value = 4.113371647056
$ws.Cells.Item(11, 19).Value = 0.09673194847494146
$ws.Cells.Item(11, 20).Value = 0.1045468662788126
